$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The image filename for the "Ironman 1/12 Plastic Model Kit FondJoy" row
# (row 7) was renamed from "fj ironman.jpg" to "ironman fj.jpg".
$ws.Range("D7").Value = "ironman fj.jpg"

# Reflect the active selection left behind after the edit.
$ws.Range("D7").Select()
